$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final Owl_Species_Table data (header + 11 data rows = 12 rows total)
$data = @(
    @("Owl_Species_ID", "Owl_Species_Common_Name"),
    @("Barn", "Barn Owl"),
    @("FerPy", "Ferruginous Pygmy-owl"),
    @("Fulvous", "Guatamalan Barred Owl"),
    @("GrHor", "Great Horned Owl"),
    @("Mottd", "Mottled Owl"),
    @("NoID", "Unidentified"),
    @("None", "None"),
    @("PacSc", "Pacific Screech-owl"),
    @("Specd", "Spectacled Owl"),
    @("Styg", "Stygian Owl"),
    @("Whisk", "Whiskered Screech-owl")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Update the defined name range to cover the new extent A1:B12
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Owl_Species_Table") {
        $n.RefersTo = "='Owl_Species_Table'!`$A`$1:`$B`$12"
    }
}
